$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "StatQuery" text (replaces the old combined stat query used in column C
# of every tab-definition row) - same new text goes into C2, C3 and C4.
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['English Setter']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Update the window zoom level and active selection to match the saved view.
$excel.ActiveWindow.Zoom = 85
$ws.Range("B4").Select() | Out-Null
